# New crime data collected - weekly refresh for 78th Precinct CompStat report.
# Updates the report header (volume number, week-covering dates) and the
# Week-to-Date / 28-Day / Year-to-Date / comparison figures for rows 15-29
# of the Crime Complaints table to reflect the newly collected week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8 ---
$ws.Cells.Item(8, 1).Value = 'Volume 30   Number  35'  # A8

# --- Row 9 ---
$ws.Cells.Item(9, 3).Value = 'Report Covering the Week  8/28/2023  Through  9/3/2023'  # C9

# --- Row 15 ---
$ws.Cells.Item(15, 12).Value = 0  # L15
$ws.Cells.Item(15, 14).Value = -68.181818181818  # N15

# --- Row 16 ---
$ws.Cells.Item(16, 3).Value = 5  # C16
$ws.Cells.Item(16, 4).Value = 3  # D16
$ws.Cells.Item(16, 5).Value = 66.666666666666  # E16
$ws.Cells.Item(16, 6).Value = 15  # F16
$ws.Cells.Item(16, 7).Value = 8  # G16
$ws.Cells.Item(16, 8).Value = 87.5  # H16
$ws.Cells.Item(16, 9).Value = 89  # I16
$ws.Cells.Item(16, 10).Value = 78  # J16
$ws.Cells.Item(16, 11).Value = 14.102564102564  # K16
$ws.Cells.Item(16, 12).Value = 53.448275862069  # L16
$ws.Cells.Item(16, 13).Value = -9.183673469387  # M16
$ws.Cells.Item(16, 14).Value = -83.047619047619  # N16

# --- Row 17 ---
$ws.Cells.Item(17, 3).Value = 5  # C17
$ws.Cells.Item(17, 4).Value = 2  # D17
$ws.Cells.Item(17, 5).Value = 150  # E17
$ws.Cells.Item(17, 6).Value = 19  # F17
$ws.Cells.Item(17, 7).Value = 14  # G17
$ws.Cells.Item(17, 8).Value = 35.714285714285  # H17
$ws.Cells.Item(17, 9).Value = 138  # I17
$ws.Cells.Item(17, 10).Value = 100  # J17
$ws.Cells.Item(17, 11).Value = 38  # K17
$ws.Cells.Item(17, 12).Value = 109.090909090909  # L17
$ws.Cells.Item(17, 13).Value = 318.181818181818  # M17
$ws.Cells.Item(17, 14).Value = -23.333333333333  # N17

# --- Row 18 ---
$ws.Cells.Item(18, 3).Value = 5  # C18
$ws.Cells.Item(18, 4).Value = 2  # D18
$ws.Cells.Item(18, 5).Value = 150  # E18
$ws.Cells.Item(18, 6).Value = 13  # F18
$ws.Cells.Item(18, 7).Value = 8  # G18
$ws.Cells.Item(18, 8).Value = 62.5  # H18
$ws.Cells.Item(18, 9).Value = 129  # I18
$ws.Cells.Item(18, 10).Value = 83  # J18
$ws.Cells.Item(18, 11).Value = 55.421686746988  # K18
$ws.Cells.Item(18, 12).Value = 67.532467532467  # L18
$ws.Cells.Item(18, 13).Value = 1.574803149606  # M18
$ws.Cells.Item(18, 14).Value = -75.797373358349  # N18

# --- Row 19 ---
$ws.Cells.Item(19, 4).Value = 13  # D19
$ws.Cells.Item(19, 5).Value = -15.384615384615  # E19
$ws.Cells.Item(19, 6).Value = 42  # F19
$ws.Cells.Item(19, 7).Value = 64  # G19
$ws.Cells.Item(19, 8).Value = -34.375  # H19
$ws.Cells.Item(19, 9).Value = 392  # I19
$ws.Cells.Item(19, 10).Value = 412  # J19
$ws.Cells.Item(19, 11).Value = -4.854368932038  # K19
$ws.Cells.Item(19, 12).Value = 34.246575342465  # L19
$ws.Cells.Item(19, 13).Value = 42.028985507246  # M19
$ws.Cells.Item(19, 14).Value = 16.320474777448  # N19

# --- Row 20 ---
$ws.Cells.Item(20, 4).Value = 4  # D20
$ws.Cells.Item(20, 5).Value = -75  # E20
$ws.Cells.Item(20, 6).Value = 7  # F20
$ws.Cells.Item(20, 7).Value = 14  # G20
$ws.Cells.Item(20, 8).Value = -50  # H20
$ws.Cells.Item(20, 9).Value = 84  # I20
$ws.Cells.Item(20, 10).Value = 72  # J20
$ws.Cells.Item(20, 11).Value = 16.666666666666  # K20
$ws.Cells.Item(20, 12).Value = 86.666666666666  # L20
$ws.Cells.Item(20, 13).Value = 35.483870967741  # M20
$ws.Cells.Item(20, 14).Value = -88.218793828892  # N20

# --- Row 21 ---
$ws.Cells.Item(21, 3).Value = 27  # C21
$ws.Cells.Item(21, 4).Value = 24  # D21
$ws.Cells.Item(21, 5).Value = 12.5  # E21
$ws.Cells.Item(21, 6).Value = 97  # F21
$ws.Cells.Item(21, 7).Value = 110  # G21
$ws.Cells.Item(21, 8).Value = -11.818181818181  # H21
$ws.Cells.Item(21, 9).Value = 840  # I21
$ws.Cells.Item(21, 10).Value = 756  # J21
$ws.Cells.Item(21, 11).Value = 11.111111111111  # K21
$ws.Cells.Item(21, 12).Value = 53.846153846153  # L21
$ws.Cells.Item(21, 13).Value = 39.767054908485  # M21
$ws.Cells.Item(21, 14).Value = -63.699222126188  # N21

# --- Row 22 ---
$ws.Cells.Item(22, 3).NumberFormat = "#,##0"  # C22: style 14 -> 15 (text placeholder -> numeric)
$ws.Cells.Item(22, 3).Value = 1  # C22
$ws.Cells.Item(22, 6).Value = 2  # F22
$ws.Cells.Item(22, 8).Value = 100  # H22
$ws.Cells.Item(22, 9).Value = 17  # I22
$ws.Cells.Item(22, 11).Value = 6.25  # K22
$ws.Cells.Item(22, 12).Value = -5.555555555555  # L22
$ws.Cells.Item(22, 13).Value = -29.166666666666  # M22

# --- Row 23 ---
$ws.Cells.Item(23, 7).Value = 2  # G23
$ws.Cells.Item(23, 8).Value = 150  # H23

# --- Row 24 ---
$ws.Cells.Item(24, 3).Value = 39  # C24
$ws.Cells.Item(24, 4).Value = 31  # D24
$ws.Cells.Item(24, 5).Value = 25.806451612903  # E24
$ws.Cells.Item(24, 7).Value = 113  # G24
$ws.Cells.Item(24, 8).Value = 6.194690265486  # H24
$ws.Cells.Item(24, 9).Value = 1029  # I24
$ws.Cells.Item(24, 10).Value = 885  # J24
$ws.Cells.Item(24, 11).Value = 16.271186440678  # K24
$ws.Cells.Item(24, 12).Value = 63.074484944532  # L24
$ws.Cells.Item(24, 13).Value = 69.522240527182  # M24

# --- Row 25 ---
$ws.Cells.Item(25, 3).Value = 5  # C25
$ws.Cells.Item(25, 4).Value = 5  # D25
$ws.Cells.Item(25, 5).Value = 0  # E25
$ws.Cells.Item(25, 6).Value = 26  # F25
$ws.Cells.Item(25, 7).Value = 20  # G25
$ws.Cells.Item(25, 8).Value = 30  # H25
$ws.Cells.Item(25, 9).Value = 212  # I25
$ws.Cells.Item(25, 10).Value = 203  # J25
$ws.Cells.Item(25, 11).Value = 4.433497536945  # K25
$ws.Cells.Item(25, 12).Value = 38.562091503268  # L25
$ws.Cells.Item(25, 13).Value = 54.744525547445  # M25

# --- Row 26 ---
$ws.Cells.Item(26, 12).Value = 0  # L26

# --- Row 27 ---
$ws.Cells.Item(27, 4).Value = 2  # D27
$ws.Cells.Item(27, 6).Value = 2  # F27
$ws.Cells.Item(27, 8).Value = -60  # H27
$ws.Cells.Item(27, 10).Value = 42  # J27
$ws.Cells.Item(27, 11).Value = -23.809523809523  # K27
$ws.Cells.Item(27, 12).Value = 18.518518518518  # L27

# --- Row 28 ---
$ws.Cells.Item(28, 14).Value = -77.777777777777  # N28

# --- Row 29 ---
$ws.Cells.Item(29, 14).Value = -77.777777777777  # N29
